$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "TestCases nopCommerce" -> "TestScenario nopCommerce"
$ws.Name = "TestScenario nopCommerce"

# Update the date value in B5 (2014-02-27 -> 2024-02-27)
$ws.Range("B5").Value = "2024-02-27"

# Move the active cell/selection from D14 to D7
$ws.Range("D7").Select()
